# Automatische test-sync: 2025-08-14 20:30:50
# Append the new incoming mail-log entry to the "Logs" sheet, extend the
# conditional-formatting ranges to cover the new row, and bump the
# "Aantal" counter on the "Dashboard" sheet for the matching category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Add the new row (row 8) ---------------------------------------------
$ws.Range("A8").Value = "Nieuwe bestelling"
$ws.Range("B8").Value = "planning@testbedrijf123.nl"
$ws.Range("C8").Value = "Wil je 200 stuks M8-bouten bestellen bij onze leverancier?"
$ws.Range("D8").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E8").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F8").Value = "2025-08-14 20:30:43"
$ws.Range("G8").Value = "Nee"
$ws.Range("H8").Value = "Ja"
$ws.Range("I8").Value = "Nee"
$ws.Range("J8").Value = "Nee"

# --- Extend conditional formatting ranges from row 7 to row 8 ------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range($col + "2:" + $col + "7")
    $newRange = $ws.Range($col + "2:" + $col + "8")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard summary count -----------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 7
